# ---------------------------------------------------------------------------
# Rik Activity 2019.xlsx -- append 2019-11-xx activity log rows (663-705)
# and the new sharedStrings / table / workbook metadata that go with them.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# --- workbook-level window position (xr:revisionPtr bump happens automatically) ---
$wb.Windows.Item(1).Left = 4000
$wb.Windows.Item(1).Top = 1240

# Row 663
$ws.Range("A663").Value = 43766.59375
$ws.Range("C663").Value = "Food"
$ws.Range("D663").Value = "Chickpeas"
$ws.Range("E663").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 664
$ws.Range("A664").Value = 43766.510416666664
$ws.Range("C664").Value = "Food"
$ws.Range("D664").Value = "Couscous and chili"
$ws.Range("E664").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 665
$ws.Range("A665").Value = 43766.770833333336
$ws.Range("C665").Value = "Food"
$ws.Range("D665").Value = "Baked cod, sauteed vegetables"
$ws.Range("E665").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 666
$ws.Range("A666").Value = 43766.729166666664
$ws.Range("C666").Value = "Food"
$ws.Range("D666").Value = "Kombucha"
$ws.Range("E666").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 667
$ws.Range("A667").Value = 43767.270833333336
$ws.Range("C667").Value = "Food"
$ws.Range("D667").Value = "Latte"
$ws.Range("E667").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 668
$ws.Range("A668").Value = 43767.3125
$ws.Range("C668").Value = "Food"
$ws.Range("D668").Value = " "
$ws.Range("E668").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 669
$ws.Range("A669").Value = 43767.510416666664
$ws.Range("C669").Value = "Food"
$ws.Range("D669").Value = "Veggie bowl"
$ws.Range("E669").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 670
$ws.Range("A670").Value = 43767.729166666664
$ws.Range("C670").Value = "Food"
$ws.Range("D670").Value = "Ravioli with tomato sauce"
$ws.Range("E670").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 671
$ws.Range("A671").Value = 43767.708333333336
$ws.Range("C671").Value = "Food"
$ws.Range("D671").Value = "Pea crisps"
$ws.Range("E671").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 672
$ws.Range("A672").Value = 43766.91527777778
$ws.Range("B672").Value = 43767.259722222225
$ws.Range("C672").Value = "Sleep"
$ws.Range("E672").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 673
$ws.Range("A673").Value = 43767.92847222222
$ws.Range("B673").Value = 43768.24722222222
$ws.Range("C673").Value = "Sleep"
$ws.Range("E673").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 674
$ws.Range("A674").Value = 43768.256944444445
$ws.Range("C674").Value = "Food"
$ws.Range("D674").Value = "Latte"
$ws.Range("E674").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 675
$ws.Range("A675").Value = 43769.25
$ws.Range("C675").Value = "Food"
$ws.Range("D675").Value = "Latte"
$ws.Range("E675").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 676
$ws.Range("A676").Value = 43769.34027777778
$ws.Range("C676").Value = "Food"
$ws.Range("D676").Value = "Granola (90g) almond milk + banana"
$ws.Range("E676").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 677
$ws.Range("A677").Value = 43768.75
$ws.Range("C677").Value = "Food"
$ws.Range("D677").Value = "Indian chicken with rice"
$ws.Range("E677").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 678
$ws.Range("A678").Value = 43768.645833333336
$ws.Range("C678").Value = "Food"
$ws.Range("D678").Value = "Latte (afternoon)"
$ws.Range("E678").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 679
$ws.Range("A679").Value = 43769.416666666664
$ws.Range("C679").Value = "Food"
$ws.Range("D679").Value = "Coffee w Half and Half"
$ws.Range("E679").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 680
$ws.Range("A680").Value = 43769.80902777778
$ws.Range("C680").Value = "Food"
$ws.Range("D680").Value = "Beer"
$ws.Range("E680").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 681
$ws.Range("A681").Value = 43769.770833333336
$ws.Range("C681").Value = "Food"
$ws.Range("D681").Value = "Halibut + broccoli + bread"
$ws.Range("E681").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 682
$ws.Range("A682").Value = 43769.541666666664
$ws.Range("C682").Value = "Food"
$ws.Range("D682").Value = "Chicken tagine + rice"
$ws.Range("E682").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 683
$ws.Range("A683").Value = 43768.95940972222
$ws.Range("B683").Value = 43769.240277777775
$ws.Range("C683").Value = "Sleep"
$ws.Range("E683").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 684
$ws.Range("A684").Value = 43769.92222222222
$ws.Range("B684").Value = 43770.23888888889
$ws.Range("C684").Value = "Sleep"
$ws.Range("E684").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 685
$ws.Range("A685").Value = 43770.354166666664
$ws.Range("C685").Value = "Food"
$ws.Range("D685").Value = "Bread 132 + nut butter (60g)"
$ws.Range("E685").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 686
$ws.Range("A686").Value = 43770.6875
$ws.Range("C686").Value = "Food"
$ws.Range("D686").Value = "Pea crisps"
$ws.Range("E686").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 687
$ws.Range("A687").Value = 43771.48611111111
$ws.Range("C687").Value = "Food"
$ws.Range("D687").Value = "chickpeas"
$ws.Range("E687").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 688
$ws.Range("A688").Value = 43771.354166666664
$ws.Range("C688").Value = "Food"
$ws.Range("D688").Value = "Eggs + kimchee"
$ws.Range("E688").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 689
$ws.Range("A689").Value = 43770.91307870371
$ws.Range("B689").Value = 43771.23263888889
$ws.Range("C689").Value = "Sleep"
$ws.Range("E689").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 690
$ws.Range("A690").Value = 43770.854166666664
$ws.Range("C690").Value = "Food"
$ws.Range("D690").Value = "Bread + cheese"
$ws.Range("E690").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 691
$ws.Range("A691").Value = 43772.0125
$ws.Range("B691").Value = 43772.23125
$ws.Range("C691").Value = "Sleep"
$ws.Range("E691").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 692
$ws.Range("A692").Value = 43771.770833333336
$ws.Range("C692").Value = "Food"
$ws.Range("D692").Value = "Gnocci"
$ws.Range("E692").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 693
$ws.Range("A693").Value = 43771.541666666664
$ws.Range("C693").Value = "Food"
$ws.Range("D693").Value = "Rice + chicken curry"
$ws.Range("E693").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 694
$ws.Range("A694").Value = 43772.510416666664
$ws.Range("C694").Value = "Food"
$ws.Range("D694").Value = "Indian bean curry  + avocado + corn chips"
$ws.Range("E694").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 695
$ws.Range("A695").Value = 43772.490277777775
$ws.Range("C695").Value = "Food"
$ws.Range("D695").Value = "Pear Pumpkin Banana smoothie"
$ws.Range("E695").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 696
$ws.Range("A696").Value = 43772.625
$ws.Range("C696").Value = "Food"
$ws.Range("D696").Value = "Coffee w snickers"
$ws.Range("E696").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 697
$ws.Range("A697").Value = 43772.8125
$ws.Range("C697").Value = "Food"
$ws.Range("D697").Value = "Chicken rice + curry"
$ws.Range("E697").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 698
$ws.Range("A698").Value = 43772.729166666664
$ws.Range("C698").Value = "Food"
$ws.Range("D698").Value = "Bulletproof spring rolls"
$ws.Range("E698").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 699
$ws.Range("A699").Value = 43772.677083333336
$ws.Range("C699").Value = "Food"
$ws.Range("D699").Value = "Bulletproof pumpkin latte"
$ws.Range("E699").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 700
$ws.Range("A700").Value = 43772.90972222222
$ws.Range("B700").Value = 43773.24444444444
$ws.Range("C700").Value = "Sleep"
$ws.Range("E700").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 701
$ws.Range("A701").Value = 43773.910416666666
$ws.Range("B701").Value = 43774.23888888889
$ws.Range("C701").Value = "Sleep"
$ws.Range("E701").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 702
$ws.Range("A702").Value = 43773.270833333336
$ws.Range("C702").Value = "Food"
$ws.Range("D702").Value = "English muffin w/egg"
$ws.Range("E702").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 703
$ws.Range("A703").Value = 43773.375
$ws.Range("B703").Value = 43773.40625
$ws.Range("C703").Value = "Exercise"
$ws.Range("E703").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 704
$ws.Range("A704").Value = 43773.520833333336
$ws.Range("C704").Value = "Food"
$ws.Range("D704").Value = "Chicken rice + curry"
$ws.Range("E704").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 705
$ws.Range("A705").Value = 43773.770833333336
$ws.Range("C705").Value = "Food"
$ws.Range("D705").Value = "Mexican pork + frijoles"
$ws.Range("E705").Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# --- restore the date/time display format on the new Start/End cells ---
$ws.Range("A663:A705").NumberFormat = "m/d/yy h:mm"
$ws.Range("B672:B673").NumberFormat = "m/d/yy h:mm"
$ws.Range("B683:B684").NumberFormat = "m/d/yy h:mm"
$ws.Range("B689").NumberFormat = "m/d/yy h:mm"
$ws.Range("B691").NumberFormat = "m/d/yy h:mm"
$ws.Range("B700:B701").NumberFormat = "m/d/yy h:mm"
$ws.Range("B703").NumberFormat = "m/d/yy h:mm"

# --- move the visible viewport / selection to the new bottom of the sheet ---
$ws.Application.GoTo($ws.Range("A706"), $true)
$ws.Range("A706").Select()
